$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one data row per row 2..254 (row 1 = header).
# This edit duplicates row 137 (Feria Lagunitas de Puerto Montt, Limón,
# fecha 44370, 1a plateado) by copying it and inserting the copy right
# after it at row 138. That pushes the former rows 138..254 down to
# 139..255, growing the used range from A1:T254 to A1:T255 — matching
# the "Fruta / hortaliza, semanal" commit that added one extra weekly
# observation.
$ws.Rows.Item(137).Copy()
$ws.Rows.Item(138).Insert()
